$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.928.75"

$ws.Range("D3").Value = "1.632.81"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'214.63"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "'0.255"
$ws.Range("E8").Value = "  +0.32%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").Value = "1.858.96"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "'4.24"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").Value = "1.633.57"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'62.84"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0755"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "25.918.27"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Value = "'193.45"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").Value = "'4.38"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").Value = "'6.27"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").Value = "'142.37"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("E27").Value = "  +2.68%  "

$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "1.136.10"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").Value = "'0.804"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.47"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'99.19"
$ws.Range("E44").Value = "  -1.49%  "

$ws.Range("D45").Value = "1.768.03"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").Value = "'56.19"
$ws.Range("E47").Value = "  +2.21%  "

$ws.Range("E48").Value = "  +3.60%  "

$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "'0.415"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("D51").Value = "'7.63"
$ws.Range("E51").Value = "  +2.82%  "

